$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Zaragoza (Saragossa)" city-name entry should just read "Zaragoza"
# (cell D53 holds that shared string).
$ws.Range("D53").Value = "Zaragoza"

# Scroll the sheet down to the bottom of the data (the view used to be
# pinned at the top, now it is scrolled so row 91 is the first visible
# row) and leave the selection on the last data row's longitude cell.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 91
$win.ScrollColumn = 1
$ws.Range("F105").Select()

# The globe now requires every yearly dataset to be the same length,
# which nudges the (shared) column width slightly narrower.
$ws.Columns("A:AMK").ColumnWidth = 13.333333333333334
